$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Infra" sheet: rename CAPEX/OPEX headers, add Long rail / Short rail rows
# ---------------------------------------------------------------------------
$wsInfra = $wb.Worksheets.Item("Infra")

$wsInfra.Range("A5").Value = "Long rail"
$wsInfra.Range("B5").Value = 2000000
$wsInfra.Range("C5").Value = 60000

$wsInfra.Range("A6").Value = "Short rail"
$wsInfra.Range("B6").Value = 1000000
$wsInfra.Range("C6").Value = 60000

$wsInfra.Range("C1").Value = "OPEX (euros/km/year)"
$wsInfra.Range("B1").Value = "CAPEX (euros/km)"

$wsInfra.Columns.Item(2).ColumnWidth = 18.26
$wsInfra.Columns.Item(3).ColumnWidth = 18.666666666666668

$wsInfra.Range("B1").Select() | Out-Null

# ---------------------------------------------------------------------------
# "Global" sheet: add Grid construction allowed row
# ---------------------------------------------------------------------------
$wsGlobal = $wb.Worksheets.Item("Global")

$wsGlobal.Range("A6").Value = "Grid construction allowed"
$wsGlobal.Range("B6").Value = $true

$wsGlobal.Range("A7").Select() | Out-Null
